# Edit: Tabel_Prob_Stat.xlsx
# - Remove "Tema_6" column (old column H), shifting Proiect/Examen/Total left.
# - Populate the (new) Examen column (I) with real exam scores.
# - Fix a couple of mistyped Tema_2/Tema_3 values in row 24.
# - Recompute the Total formula (new column J) to use the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the entire "Tema_6" column (old column H). Everything to the
#    right (Proiect, Examen, Total) shifts one column to the left.
$ws.Range("H1").EntireColumn.Delete()

# 2) Fix two grades on row 24 (Tema_2 / Tema_3) that were left blank (0)
#    before and now have real values.
$ws.Range("D24").Value = 51
$ws.Range("E24").Value = 28

# 3) Fill in the "Examen" column (now column I) with the real exam scores,
#    keyed by worksheet row number.
$examenScores = @{
    2 = 31;
    3 = 25;
    4 = 17;
    5 = 7;
    6 = 12;
    7 = 20;
    8 = 21;
    9 = 1;
    10 = 1;
    11 = 0;
    12 = 14;
    13 = 21;
    14 = 25;
    15 = 0;
    16 = 19;
    17 = 13;
    18 = 18;
    19 = 20;
    20 = 0;
    21 = 21;
    22 = 0;
    23 = 19;
    24 = 10;
    25 = 11;
    26 = 0;
    27 = 20;
    28 = 9;
    29 = 14;
    30 = 7;
    31 = 19;
    32 = 0;
    33 = 22;
    34 = 21;
    35 = 17;
    36 = 14;
    37 = 20;
    38 = 10;
    39 = 0;
    40 = 20;
    41 = 0;
    42 = 16;
    43 = 0;
    44 = 0;
    45 = 18;
    46 = 22;
    47 = 30;
    48 = 0;
    49 = 14;
    50 = 28;
    51 = 0;
    52 = 21;
    53 = 27;
    54 = 0;
    55 = 26;
    56 = 28;
    57 = 20;
    58 = 4;
    59 = 20;
    60 = 15;
    61 = 21;
    62 = 22;
    63 = 24;
    64 = 28;
    65 = 5;
    66 = 0;
    67 = 10;
    68 = 24;
    69 = 0;
    70 = 0;
    71 = 6;
    72 = 23;
    73 = 35;
    74 = 9;
    75 = 4;
    76 = 39;
    77 = 9;
    78 = 19;
    79 = 22;
    80 = 16;
    81 = 22;
    82 = 10;
    83 = 0;
    84 = 8;
    85 = 9;
    86 = 14;
    87 = 8;
    88 = 8;
    89 = 27;
    90 = 15;
    91 = 16;
    92 = 33;
    93 = 37;
    94 = 15;
    95 = 23;
    96 = 21;
    97 = 8;
    98 = 9;
    99 = 14;
    100 = 19;
    101 = 28;
    102 = 38;
    103 = 10;
    104 = 26;
    105 = 0;
    106 = 8
}
foreach ($r in $examenScores.Keys) {
    $ws.Cells.Item($r, 9).Value = $examenScores[$r]
}

# 4) Recompute the "Total" column (now column J) with the updated formula
#    (Examen is now divided by 5 before being weighted).
$ws.Range("J2:J106").Formula = "=MIN(10, 1+0.5*I2/5+0.2*H2+0.06*(C2/9+D2/7+E2/7+F2/7+G2/6))"

# 5) Re-apply the AutoFilter so its range shrinks from A1:K106 to A1:J106.
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:J106").AutoFilter()

# 6) Update the (now stale) _FilterDatabase defined name to match.
$wb.Names.Item("Tabel_Prob_Stat!_FilterDatabase").RefersTo = "=Tabel_Prob_Stat!`$A`$1:`$J`$106"

# 7) Restore the cursor/selection to where the author left it.
$ws.Range("I76").Select()
